# Bugfixed the naive forecaster component module
#
# The source data window shifted forward by one period (the oldest
# date_of_forecast row dropped off the front) and the y_1_forecast (AR2)
# values were recomputed, which also pushed the start of the forecast
# column a few rows further down. The previously last row (old row 19)
# is no longer part of the series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final contents for rows 2..18 (row 1 is the unchanged header row).
# Columns: A = date_of_forecast (serial), B = y_0, C = y_0_forecast,
#          D = y_1, E = y_1_forecast
$data = @(
    @(39765, 2008, $null,                2009, $null),
    @(40130, 2009, -3.872359107260159,   2010, $null),
    @(40494, 2010, 4.530477057343663,    2011, $null),
    @(40862, 2011, 6.833902841285977,    2012, $null),
    @(41228, 2012, 4.166536506645224,    2013, 2.693188401769642),
    @(41592, 2013, 2.669880057548091,    2014, 3.947916604971446),
    @(41957, 2014, 5.50293301232252,     2015, 4.998814576944932),
    @(42321, 2015, 4.829481320500406,    2016, 4.673582741620552),
    @(42689, 2016, 5.100281927437122,    2017, 4.372458986620376),
    @(43053, 2017, 5.161358932333737,    2018, 4.927320050172312),
    @(43418, 2018, 5.902681694119694,    2019, 4.339089271348406),
    @(43783, 2019, 3.884502719230132,    2020, 3.243024666552685),
    @(44159, 2020, -3.840397826549158,   2021, 0.2915162802050064),
    @(44525, 2021, 0.4839811651348835,   2022, 3.818597641626909),
    @(44890, 2022, 2.06342951900429,     2023, 0.6985632195332103),
    @(45254, 2023, -2.156362896191677,   2024, 0.3452735157291054),
    @(45618, 2024, -0.8205034771073372,  2025, 1.5902148106679)
)

$cols = @("A", "B", "C", "D", "E")

$r = 2
foreach ($row in $data) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cell = $ws.Range($cols[$i] + $r)
        $val = $row[$i]
        if ($null -eq $val) {
            $cell.ClearContents()
        } else {
            $cell.Value = $val
        }
    }
    $r = $r + 1
}

# The series lost its final (19th) row entirely, so remove it and let
# everything below (there is nothing) shift up; dimension shrinks to E18.
$ws.Rows.Item(19).Delete()

Write-Host "Applied naive forecaster bugfix: shifted series window and recomputed y_1_forecast values."
